$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210. This shifts existing rows 210-226 down
# to 211-227, carrying their values/formatting along (matches the diff,
# which shows every row from 210 downward effectively moved down by one,
# with a brand-new record now occupying row 210).
$ws.Rows(210).Insert()

# Populate the newly inserted row 210 with the new weekly record.
$ws.Range("A210").Value2 = 10
$ws.Range("B210").Value2 = "Vega Modelo de Temuco"
$ws.Range("C210").Value2 = "La Araucanía"
$ws.Range("D210").Value2 = 44461
$ws.Range("E210").Value2 = 9
$ws.Range("F210").Value2 = 100112008
$ws.Range("G210").Value2 = "Coliflor"
$ws.Range("H210").Value2 = "Sin especificar"
$ws.Range("I210").Value2 = "Primera"
$ws.Range("J210").Value2 = 2000
$ws.Range("K210").Value2 = 800
$ws.Range("L210").Value2 = 800
$ws.Range("M210").Value2 = 800
$ws.Range("N210").Value2 = "$/unidad"
$ws.Range("O210").Value2 = "Región de O'Higgins"
$ws.Range("P210").Value2 = 800
$ws.Range("Q210").Value2 = 1
$ws.Range("R210").Value2 = "Hortaliza"
